$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.539.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "'1.918.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'245.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4804"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.2900"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "'0.06726"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "'110.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("D11").Value = "'19.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("D12").Value = "'1.915.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'0.07572"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "'5.273"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "'0.6680"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "'298.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "'30.518.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'13.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'5.583"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.77%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'0.000007580"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'2.162.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'0.9996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'6.478"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.43%  "
$ws.Range("D25").Value = "'9.490"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "'164.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("D27").Value = "'20.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.75%  "
$ws.Range("D28").Value = "'2.116"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'0.1078"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").Value = "'4.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'4.045"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'0.04999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'0.7368"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "'0.9998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'2.723"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "'0.02033"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").Value = "'2.685"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "'111.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("D42").Value = "'0.4434"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("D43").Value = "'72.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.27%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "'5.897"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D47").Value = "'49.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "'7.279"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").Value = "'9.318"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'0.1230"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'0.2547"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.87%  "
